$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new job posting row (JD_003) based on the existing "Senior Engineer" posting,
# but with Total_Years_Min_Exp=1, Total_Years_Max_Exp=4, and Work_Mode="Hybrid".
$jobDescription = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Senior Engineer"
$ws.Range("C4").Value = $jobDescription
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Hybrid"
$ws.Range("G4").Value = "Pune, Maharashtra, India"
